$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4879611385300813
$ws.Range("C2").Value = 0.9856673031158567
$ws.Range("D2").Value = 0.5137558584591241
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=7, n_estimators=150))])"
$ws.Range("G2").Value = 0.1242467469831657
$ws.Range("H2").Value = 0.991
